# [Fonds de solidarite] Add 2020-12-31 data
# Updates nombre_aides (col C) and montant_total (col D) for the rows whose
# underlying aggregates changed with the new 2020-12-31 data refresh.
# Values are written with a leading apostrophe so Excel stores them as text
# (matching the workbook's existing inlineStr/text-typed numeric columns),
# then the cell style is reset to "Normal" so no stray number-format / style
# is left behind on cells that were otherwise unstyled.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.Value = "'" + $value
    $rng.Style = "Normal"
}

Set-TextValue "C21" "48"
Set-TextValue "D21" "130281.00"

Set-TextValue "C24" "494"
Set-TextValue "D24" "3806709.26"

Set-TextValue "C39" "578"
Set-TextValue "D39" "6190211.97"

Set-TextValue "C47" "96"
Set-TextValue "D47" "1028987.72"

Set-TextValue "C130" "1125"
Set-TextValue "D130" "9105363.67"

Set-TextValue "C214" "1015"
Set-TextValue "D214" "9920395.79"

Set-TextValue "C222" "210"
Set-TextValue "D222" "2239125.45"

Set-TextValue "C245" "475"
Set-TextValue "D245" "3465751.80"

Set-TextValue "C249" "104"
Set-TextValue "D249" "809165.53"

Set-TextValue "C253" "92"
Set-TextValue "D253" "818127.57"

Set-TextValue "C261" "1710"
Set-TextValue "D261" "11976409.68"
